# feat: add 2022-Q4 data
#
# 1. Insert a new "2022-Q4" row at the top of the data in the "总计"
#    sheet (existing rows shift down one).
# 2. Insert a brand-new "2022-Q4" worksheet (positioned right after "总计",
#    before "2022-Q2") holding the per-fund holdings detail for that
#    quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Part 1: "总计" (summary) sheet - insert new row 2 for 2022-Q4
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Rows.Item(2).Insert()

$summary.Cells.Item(2,1).Value = 1
$summary.Cells.Item(2,2).Value = "2022-Q4"
$summary.Cells.Item(2,3).Value = 12
$summary.Cells.Item(2,4).Value = 4.4

# Inserting a row copies the formatting of the row above (the bold/bordered
# header) onto every new cell; clear it off the plain data cells, then
# restore the "index column" look (bold + border) on A2 from a known-good
# sibling (A3) so it matches A3:A7.
$summary.Cells.Item(2,2).ClearFormats()
$summary.Cells.Item(2,3).ClearFormats()
$summary.Cells.Item(2,4).ClearFormats()

$summary.Cells.Item(3,1).Copy()
$summary.Cells.Item(2,1).PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Part 2: brand-new "2022-Q4" worksheet with fund-holdings detail
# ---------------------------------------------------------------------------
$q4sheet = $wb.Worksheets.Add($wb.Worksheets.Item("2022-Q2"))
$q4sheet.Name = "2022-Q4"

# Match the sheetPr block ("<outlinePr summaryBelow="1" summaryRight="1"/>")
# that every other sheet in this workbook carries.
$q4sheet.Outline.SummaryRow = 1
$q4sheet.Outline.SummaryColumn = 1

# Borrow the header-row formatting (bold + border, B1:H1) from the existing
# "2022-Q2" sheet, which shares the same 8-column layout. NOTE: re-fetch the
# sheet reference here (rather than reusing the handle captured above) - a
# worksheet object grabbed before Worksheets.Add() goes stale for
# Range/Copy purposes once the collection has changed.
$q2sheet = $wb.Worksheets.Item("2022-Q2")
$q2sheet.Range("B1:H1").Copy()
$q4sheet.Range("B1:H1").PasteSpecial(-4122)

# "2022-Q2" only has 6 data rows, but this sheet needs 12 - so the
# index-column (A) style is stamped per-row from a single source cell
# instead of bulk-copying a range that doesn't cover every destination row.
$q2sheet.Range("A2").Copy()
$q4sheet.Range("A2:A13").PasteSpecial(-4122)

# Match the rest of the workbook's page margins (a freshly added sheet
# otherwise defaults to Excel's normal template margins).
$q4sheet.PageSetup.LeftMargin = 54
$q4sheet.PageSetup.RightMargin = 54
$q4sheet.PageSetup.TopMargin = 72
$q4sheet.PageSetup.BottomMargin = 72
$q4sheet.PageSetup.HeaderMargin = 36
$q4sheet.PageSetup.FooterMargin = 36

$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($col = 0; $col -lt $headers.Length; $col++) {
    $q4sheet.Cells.Item(1, $col + 2).Value = $headers[$col]
}

# index, fund code, fund name, scale, stock position, position ratio,
# market value (亿元), position rank
$rows = @(
    @(0,  "519697", "交银优势行业混合",             "62.83", "77.22", "2.71", "1.7027", 10),
    @(1,  "005001", "交银施罗德持续成长主题混合",     "39.94", "78.78", "2.64", "1.0544", 10),
    @(2,  "506006", "汇添富科创板2年定开混合",       "17.32", "79.02", "5.54", "0.9595", 4),
    @(3,  "014611", "富国核心科技12个月持有期混合A",  "5.97",  "81.00", "5.85", "0.3492", 1),
    @(4,  "016524", "招商均衡成长混合A",             "3.47",  "80.67", "4.01", "0.1391", 3),
    @(5,  "013630", "嘉实均衡臻选一年持有期混合A",    "1.92",  "81.00", "5.85", "0.1123", 1),
    @(6,  "014612", "富国核心科技12个月持有期混合C",  "0.56",  "81.00", "5.85", "0.0328", 1),
    @(7,  "015429", "中银证券专精特新股票A",         "1.18",  "61.04", "2.18", "0.0257", 5),
    @(8,  "016525", "招商均衡成长混合C",             "0.41",  "80.67", "4.01", "0.0164", 3),
    @(9,  "013584", "招商品质领航混合C",             "0.05",  "62.41", "6.22", "0.0031", 2),
    @(10, "015430", "中银证券专精特新股票C",         "0.14",  "61.04", "2.18", "0.0031", 5),
    @(11, "013583", "招商品质领航混合A",             "0.02",  "62.41", "6.22", "0.0012", 2)
)

# A never-touched cell far outside the used range - a cheap "style-less"
# clipboard source so numeric-looking text values (fund codes / decimal
# figures stored as text, leading zeros and all) can be forced to Text
# without leaving the plain data cells with a stray custom style.
$blank = $q4sheet.Cells.Item(200, 200)
$blank.Copy()

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $data = $rows[$i]

    $q4sheet.Cells.Item($r, 1).Value = $data[0]

    $q4sheet.Cells.Item($r, 2).Value = "'" + $data[1]
    $q4sheet.Cells.Item($r, 2).PasteSpecial(-4122)

    $q4sheet.Cells.Item($r, 3).Value = $data[2]

    $q4sheet.Cells.Item($r, 4).Value = "'" + $data[3]
    $q4sheet.Cells.Item($r, 4).PasteSpecial(-4122)

    $q4sheet.Cells.Item($r, 5).Value = "'" + $data[4]
    $q4sheet.Cells.Item($r, 5).PasteSpecial(-4122)

    $q4sheet.Cells.Item($r, 6).Value = "'" + $data[5]
    $q4sheet.Cells.Item($r, 6).PasteSpecial(-4122)

    $q4sheet.Cells.Item($r, 7).Value = "'" + $data[6]
    $q4sheet.Cells.Item($r, 7).PasteSpecial(-4122)

    $q4sheet.Cells.Item($r, 8).Value = $data[7]
}

# Restore the originally-selected tab (the last sheet, "2021-Q2") so the
# freshly-added sheet doesn't steal the active/selected tab marker.
$wb.Worksheets.Item("2021-Q2").Activate()
